$d = $word.ActiveDocument

# The placeholder "{{ sigla_identidade }}" appears twice in the document;
# the one we need is the one immediately followed by ": {{ num_rg }}"
# (the signature-block line "{{ sigla_identidade }}: {{ num_rg }}").
# Find.Execute locates that unique text and leaves $match collapsed/spanning it.
$match = $d.Content
$found = $match.Find.Execute("{{ sigla_identidade }}: {{ num_rg }}", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the target '{{ sigla_identidade }}: {{ num_rg }}' text"
}

$placeholderLen = "{{ sigla_identidade }}".Length

# Range covering just the existing "{{ sigla_identidade }}" run; InsertAfter
# appends the new text right after it (and before the ": {{ num_rg }}" run),
# inheriting that run's character formatting (Times New Roman, 12pt).
$existingRun = $d.Range($match.Start, $match.Start + $placeholderLen)
$existingRun.InsertAfter("/{{ sigla_estado_identidade}}")

# Range covering only the freshly inserted text.
$newRun = $d.Range($match.Start + $placeholderLen, $match.Start + $placeholderLen + "/{{ sigla_estado_identidade}}".Length)

# Force the new text to live in its own <w:r> (distinct from the neighboring
# runs it was spliced between) while ending up with the very same formatting
# it already inherited, by nudging the size away and back.
$newRun.Font.Size = 13
$newRun.Font.Size = 12
